# Applies the "fix text for commontoolbox" edit:
#  - exclusion_rules: add new "Ignore Store Policy" column (E) with flags,
#    add two new KPI exclusion rows, and fix row 15's brand value to "PRINGLES"
#  - store_policy: add two new KPI rows ("Some KPI", "PepsiCo Segment Space to
#    Sales Index") and make it the active/selected sheet
#  - store_policy_old remains a verbatim copy (no data changes)

$wb = $excel.ActiveWorkbook

$brandList = 'BLACK COUNTRY SNACKS, AKSAM PALUSZKI, LAJKONIK PALUSZKI, ALKA ELEPHANT, RUMPLERS, TABITHA, CHEF''S LARDER, BOROMIR, COFRESH, JIFFY POP, CROCO, THE CURATORS, FRESHERS, FUDCO, GEFEN, GINNI''S, OH MY GURU!, HALDIRAMS, HALDIRAMS SNACKS, INDIE BAY SNACKS, INNATE, JACK-LNK''S, EAZY-PP-PPCRN, ZWEIFEL CRISPS, HUNKY DORYS CRISPS, LAJKONIK JUNIOR, LOVE CHIN CHIN, NISHAS SNACKS, NUTELLA, OSEM SAVOURY SNACK, OUR LITTLE REBELL!ON, EPIC, CRAWFORDS, FLIPZ, OATIS, RYMUT SNACKS, GINCO, SUNSHINE SNACKS, JAY''S, MIDLAND SNACKS, RED MILL SAVOURY SNACKS, SENSIBLE PORTIONS, VISCONTI SNACKS, WELL & TRULY SNACKS, WILD WEST, WILDING''S, BLUE DRAGON, BEPPS, BLUE DIAMOND, COFRESH SNACKS, SCHAR, OLD EL PASO, PLANTERS, LINWOODS, CYPRESSA, KOHINOOR SNACKS, KOIKEYA, PALUSZKI, LORENZ CRISPS, MCCOLGAN, ITSU, NAIRNS, NATURES STORE SNACKS, NIM''S, BAMBA SNACKS, BISSLI SNACKS, SHARWOODS, MR PORKY SNACKS, TYGRYSKI, THE REAL PORK CRACKLING CO SNACKS, THE SNAFFLING PIG CO, WHITWORTHS, YUM & YAY'

# ---------------------------------------------------------------------------
# Sheet 1: exclusion_rules
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("exclusion_rules")

# New column header
$ws1.Cells.Item(1, 5).Value = "Ignore Store Policy"

# "Ignore Store Policy" flag values for existing rows 2-22
$ignoreFlags = @{
    2  = 1
    3  = 1
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 1
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
}
foreach ($r in $ignoreFlags.Keys) {
    $ws1.Cells.Item($r, 5).Value = $ignoreFlags[$r]
}

# Fix row 15 (Brand SOS / Exclude / brand_name) value -> "PRINGLES"
$ws1.Cells.Item(15, 4).Value = "PRINGLES"

# New row 23: Some KPI / Exclude / brand_name / <brand list> / Ignore=0
$ws1.Cells.Item(23, 1).Value = "Some KPI"
$ws1.Cells.Item(23, 2).Value = "Exclude"
$ws1.Cells.Item(23, 3).Value = "brand_name"
$ws1.Cells.Item(23, 4).Value = $brandList
$ws1.Cells.Item(23, 5).Value = 0

# New row 24: PepsiCo Segment Space to Sales Index / Exclude / brand_name / DORITOS / Ignore=0
$ws1.Cells.Item(24, 1).Value = "PepsiCo Segment Space to Sales Index"
$ws1.Cells.Item(24, 2).Value = "Exclude"
$ws1.Cells.Item(24, 3).Value = "brand_name"
$ws1.Cells.Item(24, 4).Value = "DORITOS"
$ws1.Cells.Item(24, 5).Value = 0

# ---------------------------------------------------------------------------
# Sheet 3: store_policy
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("store_policy")

# New row 6: Some KPI / TT
$ws3.Cells.Item(6, 1).Value = "Some KPI"
$ws3.Cells.Item(6, 2).Value = "TT"

# New row 7: PepsiCo Segment Space to Sales Index / TT
$ws3.Cells.Item(7, 1).Value = "PepsiCo Segment Space to Sales Index"
$ws3.Cells.Item(7, 2).Value = "TT"

# ---------------------------------------------------------------------------
# Selection / active sheet bookkeeping (matches the saved view state)
# ---------------------------------------------------------------------------
$null = $ws1.Range("A25").Select()
$null = $ws3.Activate()
$null = $ws3.Range("A19").Select()
